$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.484.38'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").Value = '2.949.97'
$ws.Range("E3").Value = '  -2.01%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.99'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.02'
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("D9").Value = '2.946.72'
$ws.Range("E9").Value = '  -1.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.68'
$ws.Range("E10").Value = '  -2.80%  '
$ws.Range("E11").Value = '  -4.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.465'
$ws.Range("E12").Value = '  +1.34%  '
$ws.Range("E13").Value = '  -2.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.82'
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("D16").Value = '65.569.30'
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("D17").Value = '3.441.42'
$ws.Range("E17").Value = '  -1.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.10'
$ws.Range("E18").Value = '  +1.89%  '
$ws.Range("D19").Value = '2.953.98'
$ws.Range("E19").Value = '  -2.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.83'
$ws.Range("E20").Value = '  +13.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '445.70'
$ws.Range("E21").Value = '  -2.62%  '
$ws.Range("E23").Value = '  -1.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.01'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.23'
$ws.Range("E25").Value = '  -3.69%  '
$ws.Range("E26").Value = '  -0.95%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.01'
$ws.Range("E28").Value = '  -5.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.52'
$ws.Range("E29").Value = '  +7.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.07'
$ws.Range("E30").Value = '  -0.82%  '
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("E32").Value = '  -0.80%  '
$ws.Range("E33").Value = '  +3.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.12'
$ws.Range("E34").Value = '  +0.26%  '
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.973'
$ws.Range("E36").Value = '  -1.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.71'
$ws.Range("E37").Value = '  -1.61%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '45.98'
$ws.Range("E38").Value = '  +5.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.98'
$ws.Range("E40").Value = '  -7.49%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.122'
$ws.Range("E41").Value = '  +0.66%  '
$ws.Range("B42").Value = 'TheGraph'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.303'
$ws.Range("E42").Value = '  -1.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.81'
$ws.Range("E43").Value = '  -6.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.51'
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '382.33'
$ws.Range("E45").Value = '  -1.65%  '
$ws.Range("E46").Value = '  -1.47%  '
$ws.Range("D47").Value = '2.681.36'
$ws.Range("E47").Value = '  -4.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.92'
$ws.Range("E48").Value = '  -1.57%  '
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.84'
$ws.Range("E50").Value = '  -0.23%  '
$ws.Range("E51").Value = '  +1.60%  '
